$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume } (only include keys that changed)
$updates = @{
    2  = @{ D = "26.255.02"; E = "  -0.35%  " }
    3  = @{ D = "1.591.72";  E = "  -0.07%  " }
    4  = @{ E = "  -0.04%  " }
    5  = @{ D = "212.91";    E = "  +0.60%  " }
    6  = @{ D = "0.501";     E = "  -0.44%  " }
    7  = @{ E = "  +0.00%  " }
    8  = @{ E = "  -0.27%  " }
    9  = @{ D = "0.0606";    E = "  -0.61%  " }
    10 = @{ D = "18.95";     E = "  -2.32%  " }
    11 = @{ D = "0.0850";    E = "  +0.53%  " }
    12 = @{ D = "1.816.11";  E = "  -0.05%  " }
    13 = @{ D = "1.598.97";  E = "  +0.94%  " }
    14 = @{ E = "  -1.22%  " }
    15 = @{ E = "  -2.56%  " }
    16 = @{ E = "  -1.16%  " }
    17 = @{ D = "26.247.91"; E = "  -0.43%  " }
    18 = @{ D = "0.0₃0724";  E = "  -1.16%  " }
    19 = @{ D = "215.15";    E = "  +1.18%  " }
    20 = @{ E = "  -1.54%  " }
    21 = @{ E = "  +0.01%  " }
    22 = @{ E = "  -0.48%  " }
    23 = @{ E = "  +0.62%  " }
    24 = @{ E = "  -3.77%  " }
    25 = @{ D = "144.75";    E = "  -0.06%  " }
    26 = @{ E = "  -0.02%  " }
    27 = @{ E = "  -1.41%  " }
    28 = @{ E = "  +0.21%  " }
    29 = @{ D = "15.11";     E = "  -0.62%  " }
    30 = @{ E = "  -2.24%  " }
    31 = @{ E = "  +0.44%  " }
    32 = @{ E = "  -0.79%  " }
    33 = @{ D = "1.410.73";  E = "  +5.22%  " }
    34 = @{ E = "  -0.14%  " }
    35 = @{ E = "  -0.47%  " }
    36 = @{ E = "  -1.58%  " }
    37 = @{ E = "  -4.30%  " }
    38 = @{ E = "  -1.10%  " }
    39 = @{ D = "0.821";     E = "  +0.29%  " }
    40 = @{ E = "  -0.41%  " }
    41 = @{ E = "  -0.01%  " }
    42 = @{ D = "0.960";     E = "  -7.92%  " }
    43 = @{ E = "  +1.09%  " }
    44 = @{ E = "  -0.25%  " }
    45 = @{ D = "1.728.04";  E = "  -0.10%  " }
    46 = @{ D = "60.87";     E = "  -1.22%  " }
    47 = @{ E = "  -0.88%  " }
    48 = @{ E = "  -1.26%  " }
    49 = @{ D = "0.0499";    E = "  -1.17%  " }
    50 = @{ D = "0.0952";    E = "  -3.35%  " }
}

# Price strings that parse as plain numbers would get silently coerced to a
# numeric cell by Excel's normal type-inference on `.Value` assignment (e.g.
# "212.91" -> 212.91). The source data keeps these as literal text (matching
# the original sheet, which stores every Price cell as text), so force the
# Text number format on just those cells before writing the new value. Cells
# whose new text could never parse as a number (multiple "." separators,
# subscript digits, etc.) don't need this and are left with their original
# General format.
$forceTextRows = @(5, 6, 9, 10, 11, 19, 25, 29, 39, 42, 46, 49, 50)
foreach ($r in $forceTextRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
